# "Ironed out a few bugs, created 'visible' column."
#
# 1. Rename the 'burp' sheet to 'Process' (and let the engine fix up the
#    formula reference on Output!C5 automatically).
# 2. On 'Input': insert a new column G called "Visible" (boolean TRUE for
#    every data row), pushing the old G ("Options") / H ("Errormessage")
#    columns one slot to the right.
# 3. Fix the Height value for the third object (Input!D4) from 9 to 8.
# 4. Tidy up the VLOOKUP formula on Output!C5 (drop the redundant
#    parentheses) now that it points at the renamed 'Process' sheet.
# 5. Re-point the active sheet/selection to match the saved UI state.

$wb = $excel.ActiveWorkbook

$wsInput   = $wb.Worksheets.Item("Input")
$wsProcess = $wb.Worksheets.Item("burp")
$wsOutput  = $wb.Worksheets.Item("Output")

# --- rename 'burp' -> 'Process' -------------------------------------------
$wsProcess.Name = "Process"

# --- Input: insert the new 'Visible' column --------------------------------
$wsInput.Columns("G:G").Insert()
$wsInput.Range("G1").Value = "Visible"
$wsInput.Range("G2").Value = $true
$wsInput.Range("G3").Value = $true
$wsInput.Range("G4").Value = $true
$wsInput.Range("G5").Value = $true

# --- Input: bug fix, object height 9 -> 8 ----------------------------------
$wsInput.Range("D4").Value = 8

# --- Output: simplify/repoint the weight formula ---------------------------
$wsOutput.Range("C5").Formula = "=VLOOKUP(Input!D5,Process!A2:B5,2,FALSE)*C3/1000"

# --- Process: give column A (material names) an explicit width -------------
$wsProcess.Columns("A:A").ColumnWidth = 11.8

# --- restore the on-screen selection / active sheet ------------------------
$wsProcess.Range("B6").Select()
$wsOutput.Range("C6").Select()

$wsInput.Activate()
$wsInput.Range("G6").Select()
